# Updates cryptos list values (price & 1h volume %) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

Set-TextValue $ws.Range("D2") "26.659.92"
$ws.Range("E2").Value = "  -7.31%  "

Set-TextValue $ws.Range("D3") "1.700.48"
$ws.Range("E3").Value = "  -5.84%  "

$ws.Range("E4").Value = "  +0.09%  "

Set-TextValue $ws.Range("D5") "219.68"
$ws.Range("E5").Value = "  -5.44%  "

Set-TextValue $ws.Range("D6") "0.5152"
$ws.Range("E6").Value = "  -12.91%  "

Set-TextValue $ws.Range("D7") "1.003"
$ws.Range("E7").Value = "  -0.03%  "

Set-TextValue $ws.Range("D8") "0.2663"
$ws.Range("E8").Value = "  -4.20%  "

Set-TextValue $ws.Range("D9") "22.24"
$ws.Range("E9").Value = "  -4.59%  "

Set-TextValue $ws.Range("D10") "0.06264"
$ws.Range("E10").Value = "  -8.33%  "

Set-TextValue $ws.Range("D11") "0.07321"
$ws.Range("E11").Value = "  -2.41%  "

Set-TextValue $ws.Range("D12") "1.701.20"
$ws.Range("E12").Value = "  -5.78%  "

Set-TextValue $ws.Range("D13") "4.509"
$ws.Range("E13").Value = "  -5.45%  "

Set-TextValue $ws.Range("D14") "0.5857"
$ws.Range("E14").Value = "  -6.10%  "

Set-TextValue $ws.Range("D15") "1.930.86"
$ws.Range("E15").Value = "  -5.88%  "

Set-TextValue $ws.Range("D16") "0.000008413"
$ws.Range("E16").Value = "  -9.36%  "

Set-TextValue $ws.Range("D17") "65.64"
$ws.Range("E17").Value = "  -13.39%  "

Set-TextValue $ws.Range("D18") "26.711.27"
$ws.Range("E18").Value = "  -7.00%  "

Set-TextValue $ws.Range("D19") "5.064"
$ws.Range("E19").Value = "  -7.61%  "

Set-TextValue $ws.Range("D21") "10.90"
$ws.Range("E21").Value = "  -5.17%  "

Set-TextValue $ws.Range("D22") "188.22"
$ws.Range("E22").Value = "  -11.00%  "

Set-TextValue $ws.Range("D23") "6.292"
$ws.Range("E23").Value = "  -8.04%  "

Set-TextValue $ws.Range("D24") "1.004"
$ws.Range("E24").Value = "  +0.05%  "

Set-TextValue $ws.Range("D25") "145.30"
$ws.Range("E25").Value = "  -5.80%  "

Set-TextValue $ws.Range("D26") "7.640"
$ws.Range("E26").Value = "  -3.03%  "

Set-TextValue $ws.Range("D27") "0.1152"
$ws.Range("E27").Value = "  -9.34%  "

Set-TextValue $ws.Range("D28") "15.80"
$ws.Range("E28").Value = "  -3.93%  "

Set-TextValue $ws.Range("D29") "1.307"
$ws.Range("E29").Value = "  -8.56%  "

Set-TextValue $ws.Range("D30") "0.05742"
$ws.Range("E30").Value = "  -7.27%  "

Set-TextValue $ws.Range("D31") "1.337"
$ws.Range("E31").Value = "  -6.15%  "

Set-TextValue $ws.Range("D32") "3.525"
$ws.Range("E32").Value = "  -6.28%  "

Set-TextValue $ws.Range("D33") "3.510"
$ws.Range("E33").Value = "  -7.31%  "

Set-TextValue $ws.Range("D34") "1.669"
$ws.Range("E34").Value = "  -3.77%  "

Set-TextValue $ws.Range("D35") "1.025"
$ws.Range("E35").Value = "  -3.73%  "

Set-TextValue $ws.Range("D36") "0.6035"
$ws.Range("E36").Value = "  -6.30%  "

Set-TextValue $ws.Range("D37") "2.375"
$ws.Range("E37").Value = "  -4.83%  "

Set-TextValue $ws.Range("D38") "2.682"
$ws.Range("E38").Value = "  -1.30%  "

Set-TextValue $ws.Range("D39") "1.094.82"
$ws.Range("E39").Value = "  -4.40%  "

Set-TextValue $ws.Range("D40") "0.01603"
$ws.Range("E40").Value = "  -6.14%  "

Set-TextValue $ws.Range("D41") "0.8696"
$ws.Range("E41").Value = "  -1.48%  "

Set-TextValue $ws.Range("D42") "5.907"
$ws.Range("E42").Value = "  -10.49%  "

$ws.Range("E43").Value = "  -0.45%  "

Set-TextValue $ws.Range("D44") "98.82"
$ws.Range("E44").Value = "  -1.34%  "

Set-TextValue $ws.Range("D45") "1.858.85"
$ws.Range("E45").Value = "  -5.23%  "

Set-TextValue $ws.Range("D48") "8.203"
$ws.Range("E48").Value = "  -1.87%  "

Set-TextValue $ws.Range("D49") "1.003"
$ws.Range("E49").Value = "  -0.23%  "

Set-TextValue $ws.Range("D50") "0.05248"
$ws.Range("E50").Value = "  -4.14%  "

Set-TextValue $ws.Range("D51") "0.4321"
$ws.Range("E51").Value = "  -3.70%  "

# Row 46/47: BabyDogeCoin moves up to rank 44, Aave drops to rank 45 (values refreshed)
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D46") "0.00000000108"
$ws.Range("E46").Value = "  -2.81%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D47") "56.97"
$ws.Range("E47").Value = "  -5.81%  "
